$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.956.02'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.834.39'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.59'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6899'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07705'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3053'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -2.59%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.54'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.81%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07808'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.850.20'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.080'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.45%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6812'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.454'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008346'
$ws.Range('D17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '28.963.06'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '243.59'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -4.09%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.083.83'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.486'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '163.31'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1470'
$ws.Range('D26').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.22'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.23%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.553'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +3.40%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.218'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.167'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.172'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05114'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7732'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.24%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.851'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01850'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.242.24'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.42%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9416'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +5.43%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '108.39'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9996'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.715'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -5.13%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '9.578'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -4.93%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.984.53'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5174'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '64.70'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -8.98%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.749'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.76%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4201'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.36%  '
$ws.Range('E51').Style = "Normal"
